$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Row 2: A 18589-2022
$ws.Range("A2").Value = 'A 18589-2022'
$ws.Range("B2").Value = 44687
$ws.Range("C2").Value = 46077
$ws.Range("D2").Value = 'SKÅNE LÄN'
$ws.Range("E2").Value = 'SIMRISHAMN'
$ws.Range("F2").Value = ''
$ws.Range("G2").Value = 2.5
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 2
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 5
$ws.Range("P2").Value = 3
$ws.Range("Q2").Value = 5
$ws.Range("R2").Value = 'Skogsalm`r`nAsk`r`nBacksippa`r`nSminkrot`r`nÅkerkulla'
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/artfynd/A 18589-2022 artfynd.xlsx", "A 18589-2022")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/kartor/A 18589-2022 karta.png", "A 18589-2022")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/klagomål/A 18589-2022 FSC-klagomål.docx", "A 18589-2022")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/klagomålsmail/A 18589-2022 FSC-klagomål mail.docx", "A 18589-2022")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/tillsyn/A 18589-2022 tillsynsbegäran.docx", "A 18589-2022")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/tillsynsmail/A 18589-2022 tillsynsbegäran mail.docx", "A 18589-2022")'

# Row 3: A 45-2022
$ws.Range("A3").Value = 'A 45-2022'
$ws.Range("B3").Value = 44564
$ws.Range("C3").Value = 46077
$ws.Range("D3").Value = 'SKÅNE LÄN'
$ws.Range("E3").Value = 'SIMRISHAMN'
$ws.Range("F3").Value = 'Övriga Aktiebolag'
$ws.Range("G3").Value = 7.2
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 4
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 4
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 5
$ws.Range("R3").Value = 'Bokvårtlav`r`nLunglav`r`nMjukdån`r`nSmåjungfrukam`r`nStor knopplav'
$ws.Range("S3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/artfynd/A 45-2022 artfynd.xlsx", "A 45-2022")'
$ws.Range("T3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/kartor/A 45-2022 karta.png", "A 45-2022")'
$ws.Range("V3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/klagomål/A 45-2022 FSC-klagomål.docx", "A 45-2022")'
$ws.Range("W3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/klagomålsmail/A 45-2022 FSC-klagomål mail.docx", "A 45-2022")'
$ws.Range("X3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/tillsyn/A 45-2022 tillsynsbegäran.docx", "A 45-2022")'
$ws.Range("Y3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/tillsynsmail/A 45-2022 tillsynsbegäran mail.docx", "A 45-2022")'

# Row 4: A 27430-2024
$ws.Range("A4").Value = 'A 27430-2024'
$ws.Range("B4").Value = 45474
$ws.Range("C4").Value = 46077
$ws.Range("D4").Value = 'SKÅNE LÄN'
$ws.Range("E4").Value = 'SIMRISHAMN'
$ws.Range("F4").Value = ''
$ws.Range("G4").Value = 3.9
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 2
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 3
$ws.Range("R4").Value = 'Spillkråka`r`nBokoxe`r`nRostfläck'
$ws.Range("S4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/artfynd/A 27430-2024 artfynd.xlsx", "A 27430-2024")'
$ws.Range("T4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/kartor/A 27430-2024 karta.png", "A 27430-2024")'
$ws.Range("V4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/klagomål/A 27430-2024 FSC-klagomål.docx", "A 27430-2024")'
$ws.Range("W4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/klagomålsmail/A 27430-2024 FSC-klagomål mail.docx", "A 27430-2024")'
$ws.Range("X4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/tillsyn/A 27430-2024 tillsynsbegäran.docx", "A 27430-2024")'
$ws.Range("Y4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/tillsynsmail/A 27430-2024 tillsynsbegäran mail.docx", "A 27430-2024")'
$ws.Range("Z4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/fåglar/A 27430-2024 prioriterade fågelarter.docx", "A 27430-2024")'

# Row 5: A 19874-2024
$ws.Range("A5").Value = 'A 19874-2024'
$ws.Range("B5").Value = 45433
$ws.Range("C5").Value = 46077
$ws.Range("D5").Value = 'SKÅNE LÄN'
$ws.Range("E5").Value = 'SIMRISHAMN'
$ws.Range("F5").Value = ''
$ws.Range("G5").Value = 1.6
$ws.Range("H5").Value = 3
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 0
$ws.Range("Q5").Value = 3
$ws.Range("R5").Value = 'Talltita`r`nGrönsiska`r`nKungsfågel'
$ws.Range("S5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/artfynd/A 19874-2024 artfynd.xlsx", "A 19874-2024")'
$ws.Range("T5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/kartor/A 19874-2024 karta.png", "A 19874-2024")'
$ws.Range("V5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/klagomål/A 19874-2024 FSC-klagomål.docx", "A 19874-2024")'
$ws.Range("W5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/klagomålsmail/A 19874-2024 FSC-klagomål mail.docx", "A 19874-2024")'
$ws.Range("X5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/tillsyn/A 19874-2024 tillsynsbegäran.docx", "A 19874-2024")'
$ws.Range("Y5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/tillsynsmail/A 19874-2024 tillsynsbegäran mail.docx", "A 19874-2024")'
$ws.Range("Z5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/fåglar/A 19874-2024 prioriterade fågelarter.docx", "A 19874-2024")'

# Row 6: A 39121-2021
$ws.Range("A6").Value = 'A 39121-2021'
$ws.Range("B6").Value = 44412
$ws.Range("C6").Value = 46077
$ws.Range("D6").Value = 'SKÅNE LÄN'
$ws.Range("E6").Value = 'SIMRISHAMN'
$ws.Range("F6").Value = ''
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 2
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 0
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = 0
$ws.Range("Q6").Value = 2
$ws.Range("R6").Value = 'Lövgroda`r`nStörre vattensalamander'
$ws.Range("S6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/artfynd/A 39121-2021 artfynd.xlsx", "A 39121-2021")'
$ws.Range("T6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/kartor/A 39121-2021 karta.png", "A 39121-2021")'
$ws.Range("V6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/klagomål/A 39121-2021 FSC-klagomål.docx", "A 39121-2021")'
$ws.Range("W6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/klagomålsmail/A 39121-2021 FSC-klagomål mail.docx", "A 39121-2021")'
$ws.Range("X6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/tillsyn/A 39121-2021 tillsynsbegäran.docx", "A 39121-2021")'
$ws.Range("Y6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/tillsynsmail/A 39121-2021 tillsynsbegäran mail.docx", "A 39121-2021")'

# Row 7: A 20090-2024
$ws.Range("A7").Value = 'A 20090-2024'
$ws.Range("B7").Value = 45434
$ws.Range("C7").Value = 46077
$ws.Range("D7").Value = 'SKÅNE LÄN'
$ws.Range("E7").Value = 'SIMRISHAMN'
$ws.Range("F7").Value = ''
$ws.Range("G7").Value = 2.4
$ws.Range("H7").Value = 1
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 0
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = 0
$ws.Range("Q7").Value = 1
$ws.Range("R7").Value = 'Hasselmus'
$ws.Range("S7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/artfynd/A 20090-2024 artfynd.xlsx", "A 20090-2024")'
$ws.Range("T7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/kartor/A 20090-2024 karta.png", "A 20090-2024")'
$ws.Range("V7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/klagomål/A 20090-2024 FSC-klagomål.docx", "A 20090-2024")'
$ws.Range("W7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/klagomålsmail/A 20090-2024 FSC-klagomål mail.docx", "A 20090-2024")'
$ws.Range("X7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/tillsyn/A 20090-2024 tillsynsbegäran.docx", "A 20090-2024")'
$ws.Range("Y7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/tillsynsmail/A 20090-2024 tillsynsbegäran mail.docx", "A 20090-2024")'

# Row 8: A 7082-2024
$ws.Range("A8").Value = 'A 7082-2024'
$ws.Range("B8").Value = 45343
$ws.Range("C8").Value = 46077
$ws.Range("D8").Value = 'SKÅNE LÄN'
$ws.Range("E8").Value = 'SIMRISHAMN'
$ws.Range("F8").Value = 'Övriga Aktiebolag'
$ws.Range("G8").Value = 22.7
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 0
$ws.Range("N8").Value = 0
$ws.Range("O8").Value = 1
$ws.Range("P8").Value = 1
$ws.Range("Q8").Value = 1
$ws.Range("R8").Value = 'Småvänderot'
$ws.Range("S8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/artfynd/A 7082-2024 artfynd.xlsx", "A 7082-2024")'
$ws.Range("T8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/kartor/A 7082-2024 karta.png", "A 7082-2024")'
$ws.Range("V8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/klagomål/A 7082-2024 FSC-klagomål.docx", "A 7082-2024")'
$ws.Range("W8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/klagomålsmail/A 7082-2024 FSC-klagomål mail.docx", "A 7082-2024")'
$ws.Range("X8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/tillsyn/A 7082-2024 tillsynsbegäran.docx", "A 7082-2024")'
$ws.Range("Y8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/tillsynsmail/A 7082-2024 tillsynsbegäran mail.docx", "A 7082-2024")'

# Row 9: A 45802-2022
$ws.Range("A9").Value = 'A 45802-2022'
$ws.Range("B9").Value = 44844
$ws.Range("C9").Value = 46077
$ws.Range("D9").Value = 'SKÅNE LÄN'
$ws.Range("E9").Value = 'SIMRISHAMN'
$ws.Range("F9").Value = ''
$ws.Range("G9").Value = 1.2
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0
$ws.Range("N9").Value = 0
$ws.Range("O9").Value = 1
$ws.Range("P9").Value = 1
$ws.Range("Q9").Value = 1
$ws.Range("R9").Value = 'Hartsticka'
$ws.Range("S9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/artfynd/A 45802-2022 artfynd.xlsx", "A 45802-2022")'
$ws.Range("T9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/kartor/A 45802-2022 karta.png", "A 45802-2022")'
$ws.Range("V9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/klagomål/A 45802-2022 FSC-klagomål.docx", "A 45802-2022")'
$ws.Range("W9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/klagomålsmail/A 45802-2022 FSC-klagomål mail.docx", "A 45802-2022")'
$ws.Range("X9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/tillsyn/A 45802-2022 tillsynsbegäran.docx", "A 45802-2022")'
$ws.Range("Y9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/tillsynsmail/A 45802-2022 tillsynsbegäran mail.docx", "A 45802-2022")'

# Row 10: A 951-2024
$ws.Range("A10").Value = 'A 951-2024'
$ws.Range("B10").Value = 45301
$ws.Range("C10").Value = 46077
$ws.Range("D10").Value = 'SKÅNE LÄN'
$ws.Range("E10").Value = 'SIMRISHAMN'
$ws.Range("F10").Value = 'Övriga Aktiebolag'
$ws.Range("G10").Value = 28.1
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 0
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = 0
$ws.Range("Q10").Value = 1
$ws.Range("R10").Value = 'Dvärgpipistrell'
$ws.Range("S10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/artfynd/A 951-2024 artfynd.xlsx", "A 951-2024")'
$ws.Range("T10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/kartor/A 951-2024 karta.png", "A 951-2024")'
$ws.Range("V10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/klagomål/A 951-2024 FSC-klagomål.docx", "A 951-2024")'
$ws.Range("W10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/klagomålsmail/A 951-2024 FSC-klagomål mail.docx", "A 951-2024")'
$ws.Range("X10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/tillsyn/A 951-2024 tillsynsbegäran.docx", "A 951-2024")'
$ws.Range("Y10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/tillsynsmail/A 951-2024 tillsynsbegäran mail.docx", "A 951-2024")'

# Row 11: A 2644-2026
$ws.Range("A11").Value = 'A 2644-2026'
$ws.Range("B11").Value = 46037.62291666667
$ws.Range("C11").Value = 46077
$ws.Range("D11").Value = 'SKÅNE LÄN'
$ws.Range("E11").Value = 'SIMRISHAMN'
$ws.Range("F11").Value = ''
$ws.Range("G11").Value = 0.7
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 0
$ws.Range("N11").Value = 0
$ws.Range("O11").Value = 1
$ws.Range("P11").Value = 1
$ws.Range("Q11").Value = 1
$ws.Range("R11").Value = 'Flikbålmossa'
$ws.Range("S11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/artfynd/A 2644-2026 artfynd.xlsx", "A 2644-2026")'
$ws.Range("T11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/kartor/A 2644-2026 karta.png", "A 2644-2026")'
$ws.Range("V11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/klagomål/A 2644-2026 FSC-klagomål.docx", "A 2644-2026")'
$ws.Range("W11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/klagomålsmail/A 2644-2026 FSC-klagomål mail.docx", "A 2644-2026")'
$ws.Range("X11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/tillsyn/A 2644-2026 tillsynsbegäran.docx", "A 2644-2026")'
$ws.Range("Y11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/tillsynsmail/A 2644-2026 tillsynsbegäran mail.docx", "A 2644-2026")'

# Row 12: A 3401-2024
$ws.Range("A12").Value = 'A 3401-2024'
$ws.Range("B12").Value = 45318
$ws.Range("C12").Value = 46077
$ws.Range("D12").Value = 'SKÅNE LÄN'
$ws.Range("E12").Value = 'SIMRISHAMN'
$ws.Range("F12").Value = ''
$ws.Range("G12").Value = 6.8
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 1
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 0
$ws.Range("N12").Value = 0
$ws.Range("O12").Value = 1
$ws.Range("P12").Value = 0
$ws.Range("Q12").Value = 1
$ws.Range("R12").Value = 'Borsttåg'
$ws.Range("S12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/artfynd/A 3401-2024 artfynd.xlsx", "A 3401-2024")'
$ws.Range("T12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/kartor/A 3401-2024 karta.png", "A 3401-2024")'
$ws.Range("V12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/klagomål/A 3401-2024 FSC-klagomål.docx", "A 3401-2024")'
$ws.Range("W12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/klagomålsmail/A 3401-2024 FSC-klagomål mail.docx", "A 3401-2024")'
$ws.Range("X12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/tillsyn/A 3401-2024 tillsynsbegäran.docx", "A 3401-2024")'
$ws.Range("Y12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/tillsynsmail/A 3401-2024 tillsynsbegäran mail.docx", "A 3401-2024")'

# Row 13: A 13223-2021
$ws.Range("A13").Value = 'A 13223-2021'
$ws.Range("B13").Value = 44272.58836805556
$ws.Range("C13").Value = 46077
$ws.Range("D13").Value = 'SKÅNE LÄN'
$ws.Range("E13").Value = 'SIMRISHAMN'
$ws.Range("F13").Value = ''
$ws.Range("G13").Value = 0.6
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 0
$ws.Range("N13").Value = 0
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = 0
$ws.Range("Q13").Value = 0
$ws.Range("R13").Value = ''

# Row 14: A 7258-2022
$ws.Range("A14").Value = 'A 7258-2022'
$ws.Range("B14").Value = 44606
$ws.Range("C14").Value = 46077
$ws.Range("D14").Value = 'SKÅNE LÄN'
$ws.Range("E14").Value = 'SIMRISHAMN'
$ws.Range("F14").Value = ''
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = 0
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = 0
$ws.Range("Q14").Value = 0
$ws.Range("R14").Value = ''

# Row 15: A 25881-2021
$ws.Range("A15").Value = 'A 25881-2021'
$ws.Range("B15").Value = 44344.37106481481
$ws.Range("C15").Value = 46077
$ws.Range("D15").Value = 'SKÅNE LÄN'
$ws.Range("E15").Value = 'SIMRISHAMN'
$ws.Range("F15").Value = 'Övriga Aktiebolag'
$ws.Range("G15").Value = 3.7
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = 0
$ws.Range("O15").Value = 0
$ws.Range("P15").Value = 0
$ws.Range("Q15").Value = 0
$ws.Range("R15").Value = ''

# Row 16: A 25834-2022
$ws.Range("A16").Value = 'A 25834-2022'
$ws.Range("B16").Value = 44733.58805555556
$ws.Range("C16").Value = 46077
$ws.Range("D16").Value = 'SKÅNE LÄN'
$ws.Range("E16").Value = 'SIMRISHAMN'
$ws.Range("F16").Value = ''
$ws.Range("G16").Value = 0.6
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = 0
$ws.Range("N16").Value = 0
$ws.Range("O16").Value = 0
$ws.Range("P16").Value = 0
$ws.Range("Q16").Value = 0
$ws.Range("R16").Value = ''

# Row 17: A 25822-2022
$ws.Range("A17").Value = 'A 25822-2022'
$ws.Range("B17").Value = 44733.56212962963
$ws.Range("C17").Value = 46077
$ws.Range("D17").Value = 'SKÅNE LÄN'
$ws.Range("E17").Value = 'SIMRISHAMN'
$ws.Range("F17").Value = ''
$ws.Range("G17").Value = 0.9
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 0
$ws.Range("N17").Value = 0
$ws.Range("O17").Value = 0
$ws.Range("P17").Value = 0
$ws.Range("Q17").Value = 0
$ws.Range("R17").Value = ''

# Row 18: A 17110-2025
$ws.Range("A18").Value = 'A 17110-2025'
$ws.Range("B18").Value = 45755
$ws.Range("C18").Value = 46077
$ws.Range("D18").Value = 'SKÅNE LÄN'
$ws.Range("E18").Value = 'SIMRISHAMN'
$ws.Range("F18").Value = ''
$ws.Range("G18").Value = 0.3
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = 0
$ws.Range("N18").Value = 0
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = 0
$ws.Range("Q18").Value = 0
$ws.Range("R18").Value = ''

# Row 19: A 5028-2024
$ws.Range("A19").Value = 'A 5028-2024'
$ws.Range("B19").Value = 45329
$ws.Range("C19").Value = 46077
$ws.Range("D19").Value = 'SKÅNE LÄN'
$ws.Range("E19").Value = 'SIMRISHAMN'
$ws.Range("F19").Value = 'Övriga Aktiebolag'
$ws.Range("G19").Value = 4.5
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 0
$ws.Range("N19").Value = 0
$ws.Range("O19").Value = 0
$ws.Range("P19").Value = 0
$ws.Range("Q19").Value = 0
$ws.Range("R19").Value = ''

# Row 20: A 63456-2023
$ws.Range("A20").Value = 'A 63456-2023'
$ws.Range("B20").Value = 45273
$ws.Range("C20").Value = 46077
$ws.Range("D20").Value = 'SKÅNE LÄN'
$ws.Range("E20").Value = 'SIMRISHAMN'
$ws.Range("F20").Value = ''
$ws.Range("G20").Value = 0.5
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = 0
$ws.Range("O20").Value = 0
$ws.Range("P20").Value = 0
$ws.Range("Q20").Value = 0
$ws.Range("R20").Value = ''

# Row 21: A 14000-2025
$ws.Range("A21").Value = 'A 14000-2025'
$ws.Range("B21").Value = 45740.26850694444
$ws.Range("C21").Value = 46077
$ws.Range("D21").Value = 'SKÅNE LÄN'
$ws.Range("E21").Value = 'SIMRISHAMN'
$ws.Range("F21").Value = ''
$ws.Range("G21").Value = 0.5
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = 0
$ws.Range("N21").Value = 0
$ws.Range("O21").Value = 0
$ws.Range("P21").Value = 0
$ws.Range("Q21").Value = 0
$ws.Range("R21").Value = ''

# Row 22: A 7072-2024
$ws.Range("A22").Value = 'A 7072-2024'
$ws.Range("B22").Value = 45343
$ws.Range("C22").Value = 46077
$ws.Range("D22").Value = 'SKÅNE LÄN'
$ws.Range("E22").Value = 'SIMRISHAMN'
$ws.Range("F22").Value = ''
$ws.Range("G22").Value = 5.1
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 0
$ws.Range("N22").Value = 0
$ws.Range("O22").Value = 0
$ws.Range("P22").Value = 0
$ws.Range("Q22").Value = 0
$ws.Range("R22").Value = ''

# Row 23: A 18682-2023
$ws.Range("A23").Value = 'A 18682-2023'
$ws.Range("B23").Value = 45043.60021990741
$ws.Range("C23").Value = 46077
$ws.Range("D23").Value = 'SKÅNE LÄN'
$ws.Range("E23").Value = 'SIMRISHAMN'
$ws.Range("F23").Value = ''
$ws.Range("G23").Value = 2.6
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 0
$ws.Range("N23").Value = 0
$ws.Range("O23").Value = 0
$ws.Range("P23").Value = 0
$ws.Range("Q23").Value = 0
$ws.Range("R23").Value = ''

# Row 24: A 3367-2024
$ws.Range("A24").Value = 'A 3367-2024'
$ws.Range("B24").Value = 45317
$ws.Range("C24").Value = 46077
$ws.Range("D24").Value = 'SKÅNE LÄN'
$ws.Range("E24").Value = 'SIMRISHAMN'
$ws.Range("F24").Value = ''
$ws.Range("G24").Value = 0.5
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = 0
$ws.Range("N24").Value = 0
$ws.Range("O24").Value = 0
$ws.Range("P24").Value = 0
$ws.Range("Q24").Value = 0
$ws.Range("R24").Value = ''

# Row 25: A 14994-2022
$ws.Range("A25").Value = 'A 14994-2022'
$ws.Range("B25").Value = 44657.53449074074
$ws.Range("C25").Value = 46077
$ws.Range("D25").Value = 'SKÅNE LÄN'
$ws.Range("E25").Value = 'SIMRISHAMN'
$ws.Range("F25").Value = 'Övriga Aktiebolag'
$ws.Range("G25").Value = 4.4
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 0
$ws.Range("N25").Value = 0
$ws.Range("O25").Value = 0
$ws.Range("P25").Value = 0
$ws.Range("Q25").Value = 0
$ws.Range("R25").Value = ''

# Row 26: A 4193-2024
$ws.Range("A26").Value = 'A 4193-2024'
$ws.Range("B26").Value = 45324
$ws.Range("C26").Value = 46077
$ws.Range("D26").Value = 'SKÅNE LÄN'
$ws.Range("E26").Value = 'SIMRISHAMN'
$ws.Range("F26").Value = 'Övriga Aktiebolag'
$ws.Range("G26").Value = 1.9
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = 0
$ws.Range("N26").Value = 0
$ws.Range("O26").Value = 0
$ws.Range("P26").Value = 0
$ws.Range("Q26").Value = 0
$ws.Range("R26").Value = ''

# Row 27: A 18856-2024
$ws.Range("A27").Value = 'A 18856-2024'
$ws.Range("B27").Value = 45426
$ws.Range("C27").Value = 46077
$ws.Range("D27").Value = 'SKÅNE LÄN'
$ws.Range("E27").Value = 'SIMRISHAMN'
$ws.Range("F27").Value = 'Övriga Aktiebolag'
$ws.Range("G27").Value = 3.3
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = 0
$ws.Range("N27").Value = 0
$ws.Range("O27").Value = 0
$ws.Range("P27").Value = 0
$ws.Range("Q27").Value = 0
$ws.Range("R27").Value = ''

# Row 28: A 7584-2024
$ws.Range("A28").Value = 'A 7584-2024'
$ws.Range("B28").Value = 45348
$ws.Range("C28").Value = 46077
$ws.Range("D28").Value = 'SKÅNE LÄN'
$ws.Range("E28").Value = 'SIMRISHAMN'
$ws.Range("F28").Value = ''
$ws.Range("G28").Value = 6.9
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 0
$ws.Range("N28").Value = 0
$ws.Range("O28").Value = 0
$ws.Range("P28").Value = 0
$ws.Range("Q28").Value = 0
$ws.Range("R28").Value = ''

# Row 29: A 35692-2024
$ws.Range("A29").Value = 'A 35692-2024'
$ws.Range("B29").Value = 45532
$ws.Range("C29").Value = 46077
$ws.Range("D29").Value = 'SKÅNE LÄN'
$ws.Range("E29").Value = 'SIMRISHAMN'
$ws.Range("F29").Value = ''
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = 0
$ws.Range("N29").Value = 0
$ws.Range("O29").Value = 0
$ws.Range("P29").Value = 0
$ws.Range("Q29").Value = 0
$ws.Range("R29").Value = ''

# Row 30: A 3402-2024
$ws.Range("A30").Value = 'A 3402-2024'
$ws.Range("B30").Value = 45318
$ws.Range("C30").Value = 46077
$ws.Range("D30").Value = 'SKÅNE LÄN'
$ws.Range("E30").Value = 'SIMRISHAMN'
$ws.Range("F30").Value = ''
$ws.Range("G30").Value = 0.5
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = 0
$ws.Range("N30").Value = 0
$ws.Range("O30").Value = 0
$ws.Range("P30").Value = 0
$ws.Range("Q30").Value = 0
$ws.Range("R30").Value = ''

# Row 31: A 19439-2024
$ws.Range("A31").Value = 'A 19439-2024'
$ws.Range("B31").Value = 45429
$ws.Range("C31").Value = 46077
$ws.Range("D31").Value = 'SKÅNE LÄN'
$ws.Range("E31").Value = 'SIMRISHAMN'
$ws.Range("F31").Value = ''
$ws.Range("G31").Value = 3.1
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = 0
$ws.Range("N31").Value = 0
$ws.Range("O31").Value = 0
$ws.Range("P31").Value = 0
$ws.Range("Q31").Value = 0
$ws.Range("R31").Value = ''

# Row 32: A 38846-2024
$ws.Range("A32").Value = 'A 38846-2024'
$ws.Range("B32").Value = 45547.60444444444
$ws.Range("C32").Value = 46077
$ws.Range("D32").Value = 'SKÅNE LÄN'
$ws.Range("E32").Value = 'SIMRISHAMN'
$ws.Range("F32").Value = ''
$ws.Range("G32").Value = 1.4
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = 0
$ws.Range("N32").Value = 0
$ws.Range("O32").Value = 0
$ws.Range("P32").Value = 0
$ws.Range("Q32").Value = 0
$ws.Range("R32").Value = ''

# Row 33: A 3392-2024
$ws.Range("A33").Value = 'A 3392-2024'
$ws.Range("B33").Value = 45318
$ws.Range("C33").Value = 46077
$ws.Range("D33").Value = 'SKÅNE LÄN'
$ws.Range("E33").Value = 'SIMRISHAMN'
$ws.Range("F33").Value = ''
$ws.Range("G33").Value = 2.7
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = 0
$ws.Range("N33").Value = 0
$ws.Range("O33").Value = 0
$ws.Range("P33").Value = 0
$ws.Range("Q33").Value = 0
$ws.Range("R33").Value = ''

# Row 34: A 19435-2024
$ws.Range("A34").Value = 'A 19435-2024'
$ws.Range("B34").Value = 45429
$ws.Range("C34").Value = 46077
$ws.Range("D34").Value = 'SKÅNE LÄN'
$ws.Range("E34").Value = 'SIMRISHAMN'
$ws.Range("F34").Value = ''
$ws.Range("G34").Value = 3.7
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = 0
$ws.Range("N34").Value = 0
$ws.Range("O34").Value = 0
$ws.Range("P34").Value = 0
$ws.Range("Q34").Value = 0
$ws.Range("R34").Value = ''

# Row 35: A 62413-2022
$ws.Range("A35").Value = 'A 62413-2022'
$ws.Range("B35").Value = 44923
$ws.Range("C35").Value = 46077
$ws.Range("D35").Value = 'SKÅNE LÄN'
$ws.Range("E35").Value = 'SIMRISHAMN'
$ws.Range("F35").Value = ''
$ws.Range("G35").Value = 2
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = 0
$ws.Range("N35").Value = 0
$ws.Range("O35").Value = 0
$ws.Range("P35").Value = 0
$ws.Range("Q35").Value = 0
$ws.Range("R35").Value = ''

# Row 36: A 14104-2024
$ws.Range("A36").Value = 'A 14104-2024'
$ws.Range("B36").Value = 45392.61707175926
$ws.Range("C36").Value = 46077
$ws.Range("D36").Value = 'SKÅNE LÄN'
$ws.Range("E36").Value = 'SIMRISHAMN'
$ws.Range("F36").Value = ''
$ws.Range("G36").Value = 5.3
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = 0
$ws.Range("N36").Value = 0
$ws.Range("O36").Value = 0
$ws.Range("P36").Value = 0
$ws.Range("Q36").Value = 0
$ws.Range("R36").Value = ''

# Row 37: A 14106-2024
$ws.Range("A37").Value = 'A 14106-2024'
$ws.Range("B37").Value = 45392
$ws.Range("C37").Value = 46077
$ws.Range("D37").Value = 'SKÅNE LÄN'
$ws.Range("E37").Value = 'SIMRISHAMN'
$ws.Range("F37").Value = ''
$ws.Range("G37").Value = 0.6
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = 0
$ws.Range("N37").Value = 0
$ws.Range("O37").Value = 0
$ws.Range("P37").Value = 0
$ws.Range("Q37").Value = 0
$ws.Range("R37").Value = ''

# Row 38: A 63464-2023
$ws.Range("A38").Value = 'A 63464-2023'
$ws.Range("B38").Value = 45273
$ws.Range("C38").Value = 46077
$ws.Range("D38").Value = 'SKÅNE LÄN'
$ws.Range("E38").Value = 'SIMRISHAMN'
$ws.Range("F38").Value = ''
$ws.Range("G38").Value = 0.7
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = 0
$ws.Range("N38").Value = 0
$ws.Range("O38").Value = 0
$ws.Range("P38").Value = 0
$ws.Range("Q38").Value = 0
$ws.Range("R38").Value = ''

# Row 39: A 23287-2022
$ws.Range("A39").Value = 'A 23287-2022'
$ws.Range("B39").Value = 44720
$ws.Range("C39").Value = 46077
$ws.Range("D39").Value = 'SKÅNE LÄN'
$ws.Range("E39").Value = 'SIMRISHAMN'
$ws.Range("F39").Value = ''
$ws.Range("G39").Value = 3.5
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = 0
$ws.Range("N39").Value = 0
$ws.Range("O39").Value = 0
$ws.Range("P39").Value = 0
$ws.Range("Q39").Value = 0
$ws.Range("R39").Value = ''

# Row 40: A 39238-2023
$ws.Range("A40").Value = 'A 39238-2023'
$ws.Range("B40").Value = 45162
$ws.Range("C40").Value = 46077
$ws.Range("D40").Value = 'SKÅNE LÄN'
$ws.Range("E40").Value = 'SIMRISHAMN'
$ws.Range("F40").Value = ''
$ws.Range("G40").Value = 4.2
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = 0
$ws.Range("N40").Value = 0
$ws.Range("O40").Value = 0
$ws.Range("P40").Value = 0
$ws.Range("Q40").Value = 0
$ws.Range("R40").Value = ''

# Row 41: A 42374-2023
$ws.Range("A41").Value = 'A 42374-2023'
$ws.Range("B41").Value = 45180.5609375
$ws.Range("C41").Value = 46077
$ws.Range("D41").Value = 'SKÅNE LÄN'
$ws.Range("E41").Value = 'SIMRISHAMN'
$ws.Range("F41").Value = ''
$ws.Range("G41").Value = 0.5
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = 0
$ws.Range("N41").Value = 0
$ws.Range("O41").Value = 0
$ws.Range("P41").Value = 0
$ws.Range("Q41").Value = 0
$ws.Range("R41").Value = ''

# Row 42: A 6006-2024
$ws.Range("A42").Value = 'A 6006-2024'
$ws.Range("B42").Value = 45336
$ws.Range("C42").Value = 46077
$ws.Range("D42").Value = 'SKÅNE LÄN'
$ws.Range("E42").Value = 'SIMRISHAMN'
$ws.Range("F42").Value = ''
$ws.Range("G42").Value = 2.9
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = 0
$ws.Range("N42").Value = 0
$ws.Range("O42").Value = 0
$ws.Range("P42").Value = 0
$ws.Range("Q42").Value = 0
$ws.Range("R42").Value = ''

# Row 43: A 3393-2024
$ws.Range("A43").Value = 'A 3393-2024'
$ws.Range("B43").Value = 45318
$ws.Range("C43").Value = 46077
$ws.Range("D43").Value = 'SKÅNE LÄN'
$ws.Range("E43").Value = 'SIMRISHAMN'
$ws.Range("F43").Value = ''
$ws.Range("G43").Value = 1.7
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = 0
$ws.Range("N43").Value = 0
$ws.Range("O43").Value = 0
$ws.Range("P43").Value = 0
$ws.Range("Q43").Value = 0
$ws.Range("R43").Value = ''

# Row 44: A 3395-2024
$ws.Range("A44").Value = 'A 3395-2024'
$ws.Range("B44").Value = 45318
$ws.Range("C44").Value = 46077
$ws.Range("D44").Value = 'SKÅNE LÄN'
$ws.Range("E44").Value = 'SIMRISHAMN'
$ws.Range("F44").Value = ''
$ws.Range("G44").Value = 2.4
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = 0
$ws.Range("N44").Value = 0
$ws.Range("O44").Value = 0
$ws.Range("P44").Value = 0
$ws.Range("Q44").Value = 0
$ws.Range("R44").Value = ''

# Row 45: A 3398-2024
$ws.Range("A45").Value = 'A 3398-2024'
$ws.Range("B45").Value = 45318
$ws.Range("C45").Value = 46077
$ws.Range("D45").Value = 'SKÅNE LÄN'
$ws.Range("E45").Value = 'SIMRISHAMN'
$ws.Range("F45").Value = ''
$ws.Range("G45").Value = 3.4
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = 0
$ws.Range("N45").Value = 0
$ws.Range("O45").Value = 0
$ws.Range("P45").Value = 0
$ws.Range("Q45").Value = 0
$ws.Range("R45").Value = ''

# Row 46: A 50134-2024
$ws.Range("A46").Value = 'A 50134-2024'
$ws.Range("B46").Value = 45600.44069444444
$ws.Range("C46").Value = 46077
$ws.Range("D46").Value = 'SKÅNE LÄN'
$ws.Range("E46").Value = 'SIMRISHAMN'
$ws.Range("F46").Value = ''
$ws.Range("G46").Value = 1.8
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = 0
$ws.Range("N46").Value = 0
$ws.Range("O46").Value = 0
$ws.Range("P46").Value = 0
$ws.Range("Q46").Value = 0
$ws.Range("R46").Value = ''

# Row 47: A 20950-2023
$ws.Range("A47").Value = 'A 20950-2023'
$ws.Range("B47").Value = 45061
$ws.Range("C47").Value = 46077
$ws.Range("D47").Value = 'SKÅNE LÄN'
$ws.Range("E47").Value = 'SIMRISHAMN'
$ws.Range("F47").Value = ''
$ws.Range("G47").Value = 6
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = 0
$ws.Range("N47").Value = 0
$ws.Range("O47").Value = 0
$ws.Range("P47").Value = 0
$ws.Range("Q47").Value = 0
$ws.Range("R47").Value = ''

# Row 48: A 43724-2024
$ws.Range("A48").Value = 'A 43724-2024'
$ws.Range("B48").Value = 45569.63178240741
$ws.Range("C48").Value = 46077
$ws.Range("D48").Value = 'SKÅNE LÄN'
$ws.Range("E48").Value = 'SIMRISHAMN'
$ws.Range("F48").Value = ''
$ws.Range("G48").Value = 1.7
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("M48").Value = 0
$ws.Range("N48").Value = 0
$ws.Range("O48").Value = 0
$ws.Range("P48").Value = 0
$ws.Range("Q48").Value = 0
$ws.Range("R48").Value = ''

# Row 49: A 30667-2023
$ws.Range("A49").Value = 'A 30667-2023'
$ws.Range("B49").Value = 45112
$ws.Range("C49").Value = 46077
$ws.Range("D49").Value = 'SKÅNE LÄN'
$ws.Range("E49").Value = 'SIMRISHAMN'
$ws.Range("F49").Value = 'Övriga Aktiebolag'
$ws.Range("G49").Value = 1
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = 0
$ws.Range("N49").Value = 0
$ws.Range("O49").Value = 0
$ws.Range("P49").Value = 0
$ws.Range("Q49").Value = 0
$ws.Range("R49").Value = ''

# Row 50: A 5488-2024
$ws.Range("A50").Value = 'A 5488-2024'
$ws.Range("B50").Value = 45334
$ws.Range("C50").Value = 46077
$ws.Range("D50").Value = 'SKÅNE LÄN'
$ws.Range("E50").Value = 'SIMRISHAMN'
$ws.Range("F50").Value = ''
$ws.Range("G50").Value = 0.5
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = 0
$ws.Range("N50").Value = 0
$ws.Range("O50").Value = 0
$ws.Range("P50").Value = 0
$ws.Range("Q50").Value = 0
$ws.Range("R50").Value = ''

# Row 51: A 3404-2024
$ws.Range("A51").Value = 'A 3404-2024'
$ws.Range("B51").Value = 45318
$ws.Range("C51").Value = 46077
$ws.Range("D51").Value = 'SKÅNE LÄN'
$ws.Range("E51").Value = 'SIMRISHAMN'
$ws.Range("F51").Value = ''
$ws.Range("G51").Value = 1.9
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = 0
$ws.Range("N51").Value = 0
$ws.Range("O51").Value = 0
$ws.Range("P51").Value = 0
$ws.Range("Q51").Value = 0
$ws.Range("R51").Value = ''

# Row 52: A 39015-2023
$ws.Range("A52").Value = 'A 39015-2023'
$ws.Range("B52").Value = 45162
$ws.Range("C52").Value = 46077
$ws.Range("D52").Value = 'SKÅNE LÄN'
$ws.Range("E52").Value = 'SIMRISHAMN'
$ws.Range("F52").Value = ''
$ws.Range("G52").Value = 9
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("M52").Value = 0
$ws.Range("N52").Value = 0
$ws.Range("O52").Value = 0
$ws.Range("P52").Value = 0
$ws.Range("Q52").Value = 0
$ws.Range("R52").Value = ''

# Row 53: A 5224-2024
$ws.Range("A53").Value = 'A 5224-2024'
$ws.Range("B53").Value = 45330
$ws.Range("C53").Value = 46077
$ws.Range("D53").Value = 'SKÅNE LÄN'
$ws.Range("E53").Value = 'SIMRISHAMN'
$ws.Range("F53").Value = 'Övriga Aktiebolag'
$ws.Range("G53").Value = 18.4
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = 0
$ws.Range("N53").Value = 0
$ws.Range("O53").Value = 0
$ws.Range("P53").Value = 0
$ws.Range("Q53").Value = 0
$ws.Range("R53").Value = ''

# Row 54: A 15039-2024
$ws.Range("A54").Value = 'A 15039-2024'
$ws.Range("B54").Value = 45399
$ws.Range("C54").Value = 46077
$ws.Range("D54").Value = 'SKÅNE LÄN'
$ws.Range("E54").Value = 'SIMRISHAMN'
$ws.Range("F54").Value = ''
$ws.Range("G54").Value = 1.3
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = 0
$ws.Range("N54").Value = 0
$ws.Range("O54").Value = 0
$ws.Range("P54").Value = 0
$ws.Range("Q54").Value = 0
$ws.Range("R54").Value = ''

# Row 55: A 40584-2024
$ws.Range("A55").Value = 'A 40584-2024'
$ws.Range("B55").Value = 45555.74299768519
$ws.Range("C55").Value = 46077
$ws.Range("D55").Value = 'SKÅNE LÄN'
$ws.Range("E55").Value = 'SIMRISHAMN'
$ws.Range("F55").Value = ''
$ws.Range("G55").Value = 0.8
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = 0
$ws.Range("N55").Value = 0
$ws.Range("O55").Value = 0
$ws.Range("P55").Value = 0
$ws.Range("Q55").Value = 0
$ws.Range("R55").Value = ''

# Row 56: A 26471-2023
$ws.Range("A56").Value = 'A 26471-2023'
$ws.Range("B56").Value = 45092
$ws.Range("C56").Value = 46077
$ws.Range("D56").Value = 'SKÅNE LÄN'
$ws.Range("E56").Value = 'SIMRISHAMN'
$ws.Range("F56").Value = ''
$ws.Range("G56").Value = 0.4
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = 0
$ws.Range("N56").Value = 0
$ws.Range("O56").Value = 0
$ws.Range("P56").Value = 0
$ws.Range("Q56").Value = 0
$ws.Range("R56").Value = ''

# Row 57: A 56764-2024
$ws.Range("A57").Value = 'A 56764-2024'
$ws.Range("B57").Value = 45628
$ws.Range("C57").Value = 46077
$ws.Range("D57").Value = 'SKÅNE LÄN'
$ws.Range("E57").Value = 'SIMRISHAMN'
$ws.Range("F57").Value = ''
$ws.Range("G57").Value = 7.5
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = 0
$ws.Range("N57").Value = 0
$ws.Range("O57").Value = 0
$ws.Range("P57").Value = 0
$ws.Range("Q57").Value = 0
$ws.Range("R57").Value = ''

# Row 58: A 17102-2025
$ws.Range("A58").Value = 'A 17102-2025'
$ws.Range("B58").Value = 45755
$ws.Range("C58").Value = 46077
$ws.Range("D58").Value = 'SKÅNE LÄN'
$ws.Range("E58").Value = 'SIMRISHAMN'
$ws.Range("F58").Value = ''
$ws.Range("G58").Value = 0.7
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = 0
$ws.Range("N58").Value = 0
$ws.Range("O58").Value = 0
$ws.Range("P58").Value = 0
$ws.Range("Q58").Value = 0
$ws.Range("R58").Value = ''

# Row 59: A 14702-2024
$ws.Range("A59").Value = 'A 14702-2024'
$ws.Range("B59").Value = 45397.50576388889
$ws.Range("C59").Value = 46077
$ws.Range("D59").Value = 'SKÅNE LÄN'
$ws.Range("E59").Value = 'SIMRISHAMN'
$ws.Range("F59").Value = ''
$ws.Range("G59").Value = 1.1
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("M59").Value = 0
$ws.Range("N59").Value = 0
$ws.Range("O59").Value = 0
$ws.Range("P59").Value = 0
$ws.Range("Q59").Value = 0
$ws.Range("R59").Value = ''

# Row 60: A 8721-2023
$ws.Range("A60").Value = 'A 8721-2023'
$ws.Range("B60").Value = 44978
$ws.Range("C60").Value = 46077
$ws.Range("D60").Value = 'SKÅNE LÄN'
$ws.Range("E60").Value = 'SIMRISHAMN'
$ws.Range("F60").Value = ''
$ws.Range("G60").Value = 1.7
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = 0
$ws.Range("N60").Value = 0
$ws.Range("O60").Value = 0
$ws.Range("P60").Value = 0
$ws.Range("Q60").Value = 0
$ws.Range("R60").Value = ''

# Row 61: A 7716-2024
$ws.Range("A61").Value = 'A 7716-2024'
$ws.Range("B61").Value = 45349
$ws.Range("C61").Value = 46077
$ws.Range("D61").Value = 'SKÅNE LÄN'
$ws.Range("E61").Value = 'SIMRISHAMN'
$ws.Range("F61").Value = ''
$ws.Range("G61").Value = 1.8
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = 0
$ws.Range("N61").Value = 0
$ws.Range("O61").Value = 0
$ws.Range("P61").Value = 0
$ws.Range("Q61").Value = 0
$ws.Range("R61").Value = ''

# Row 62: A 30911-2024
$ws.Range("A62").Value = 'A 30911-2024'
$ws.Range("B62").Value = 45498.5925
$ws.Range("C62").Value = 46077
$ws.Range("D62").Value = 'SKÅNE LÄN'
$ws.Range("E62").Value = 'SIMRISHAMN'
$ws.Range("F62").Value = 'Övriga Aktiebolag'
$ws.Range("G62").Value = 2.6
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = 0
$ws.Range("N62").Value = 0
$ws.Range("O62").Value = 0
$ws.Range("P62").Value = 0
$ws.Range("Q62").Value = 0
$ws.Range("R62").Value = ''

# Row 63: A 39221-2023
$ws.Range("A63").Value = 'A 39221-2023'
$ws.Range("B63").Value = 45162
$ws.Range("C63").Value = 46077
$ws.Range("D63").Value = 'SKÅNE LÄN'
$ws.Range("E63").Value = 'SIMRISHAMN'
$ws.Range("F63").Value = ''
$ws.Range("G63").Value = 3.1
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = 0
$ws.Range("N63").Value = 0
$ws.Range("O63").Value = 0
$ws.Range("P63").Value = 0
$ws.Range("Q63").Value = 0
$ws.Range("R63").Value = ''

# Row 64: A 7083-2024
$ws.Range("A64").Value = 'A 7083-2024'
$ws.Range("B64").Value = 45343
$ws.Range("C64").Value = 46077
$ws.Range("D64").Value = 'SKÅNE LÄN'
$ws.Range("E64").Value = 'SIMRISHAMN'
$ws.Range("F64").Value = 'Övriga Aktiebolag'
$ws.Range("G64").Value = 10.6
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = 0
$ws.Range("N64").Value = 0
$ws.Range("O64").Value = 0
$ws.Range("P64").Value = 0
$ws.Range("Q64").Value = 0
$ws.Range("R64").Value = ''

# Row 65: A 3602-2024
$ws.Range("A65").Value = 'A 3602-2024'
$ws.Range("B65").Value = 45320
$ws.Range("C65").Value = 46077
$ws.Range("D65").Value = 'SKÅNE LÄN'
$ws.Range("E65").Value = 'SIMRISHAMN'
$ws.Range("F65").Value = ''
$ws.Range("G65").Value = 0.7
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = 0
$ws.Range("N65").Value = 0
$ws.Range("O65").Value = 0
$ws.Range("P65").Value = 0
$ws.Range("Q65").Value = 0
$ws.Range("R65").Value = ''

# Row 66: A 61876-2024
$ws.Range("A66").Value = 'A 61876-2024'
$ws.Range("B66").Value = 45652.38005787037
$ws.Range("C66").Value = 46077
$ws.Range("D66").Value = 'SKÅNE LÄN'
$ws.Range("E66").Value = 'SIMRISHAMN'
$ws.Range("F66").Value = ''
$ws.Range("G66").Value = 1.4
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = 0
$ws.Range("N66").Value = 0
$ws.Range("O66").Value = 0
$ws.Range("P66").Value = 0
$ws.Range("Q66").Value = 0
$ws.Range("R66").Value = ''

# Row 67: A 25475-2023
$ws.Range("A67").Value = 'A 25475-2023'
$ws.Range("B67").Value = 45089
$ws.Range("C67").Value = 46077
$ws.Range("D67").Value = 'SKÅNE LÄN'
$ws.Range("E67").Value = 'SIMRISHAMN'
$ws.Range("F67").Value = ''
$ws.Range("G67").Value = 1.5
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = 0
$ws.Range("N67").Value = 0
$ws.Range("O67").Value = 0
$ws.Range("P67").Value = 0
$ws.Range("Q67").Value = 0
$ws.Range("R67").Value = ''

# Row 68: A 3396-2024
$ws.Range("A68").Value = 'A 3396-2024'
$ws.Range("B68").Value = 45318
$ws.Range("C68").Value = 46077
$ws.Range("D68").Value = 'SKÅNE LÄN'
$ws.Range("E68").Value = 'SIMRISHAMN'
$ws.Range("F68").Value = ''
$ws.Range("G68").Value = 4.4
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = 0
$ws.Range("N68").Value = 0
$ws.Range("O68").Value = 0
$ws.Range("P68").Value = 0
$ws.Range("Q68").Value = 0
$ws.Range("R68").Value = ''

# Row 69: A 3405-2024
$ws.Range("A69").Value = 'A 3405-2024'
$ws.Range("B69").Value = 45318
$ws.Range("C69").Value = 46077
$ws.Range("D69").Value = 'SKÅNE LÄN'
$ws.Range("E69").Value = 'SIMRISHAMN'
$ws.Range("F69").Value = ''
$ws.Range("G69").Value = 0.9
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = 0
$ws.Range("N69").Value = 0
$ws.Range("O69").Value = 0
$ws.Range("P69").Value = 0
$ws.Range("Q69").Value = 0
$ws.Range("R69").Value = ''

# Row 70: A 39834-2021
$ws.Range("A70").Value = 'A 39834-2021'
$ws.Range("B70").Value = 44417
$ws.Range("C70").Value = 46077
$ws.Range("D70").Value = 'SKÅNE LÄN'
$ws.Range("E70").Value = 'SIMRISHAMN'
$ws.Range("F70").Value = ''
$ws.Range("G70").Value = 1.3
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = 0
$ws.Range("N70").Value = 0
$ws.Range("O70").Value = 0
$ws.Range("P70").Value = 0
$ws.Range("Q70").Value = 0
$ws.Range("R70").Value = ''

# Row 71: A 56820-2024
$ws.Range("A71").Value = 'A 56820-2024'
$ws.Range("B71").Value = 45628
$ws.Range("C71").Value = 46077
$ws.Range("D71").Value = 'SKÅNE LÄN'
$ws.Range("E71").Value = 'SIMRISHAMN'
$ws.Range("F71").Value = ''
$ws.Range("G71").Value = 1.5
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = 0
$ws.Range("N71").Value = 0
$ws.Range("O71").Value = 0
$ws.Range("P71").Value = 0
$ws.Range("Q71").Value = 0
$ws.Range("R71").Value = ''

# Row 72: A 58985-2022
$ws.Range("A72").Value = 'A 58985-2022'
$ws.Range("B72").Value = 44903.69670138889
$ws.Range("C72").Value = 46077
$ws.Range("D72").Value = 'SKÅNE LÄN'
$ws.Range("E72").Value = 'SIMRISHAMN'
$ws.Range("F72").Value = ''
$ws.Range("G72").Value = 2.6
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = 0
$ws.Range("N72").Value = 0
$ws.Range("O72").Value = 0
$ws.Range("P72").Value = 0
$ws.Range("Q72").Value = 0
$ws.Range("R72").Value = ''

# Row 73: A 62149-2022
$ws.Range("A73").Value = 'A 62149-2022'
$ws.Range("B73").Value = 44922
$ws.Range("C73").Value = 46077
$ws.Range("D73").Value = 'SKÅNE LÄN'
$ws.Range("E73").Value = 'SIMRISHAMN'
$ws.Range("F73").Value = 'Övriga Aktiebolag'
$ws.Range("G73").Value = 2.9
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = 0
$ws.Range("N73").Value = 0
$ws.Range("O73").Value = 0
$ws.Range("P73").Value = 0
$ws.Range("Q73").Value = 0
$ws.Range("R73").Value = ''

# Row 74: A 41157-2025
$ws.Range("A74").Value = 'A 41157-2025'
$ws.Range("B74").Value = 45898.52972222222
$ws.Range("C74").Value = 46077
$ws.Range("D74").Value = 'SKÅNE LÄN'
$ws.Range("E74").Value = 'SIMRISHAMN'
$ws.Range("F74").Value = ''
$ws.Range("G74").Value = 1.5
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = 0
$ws.Range("N74").Value = 0
$ws.Range("O74").Value = 0
$ws.Range("P74").Value = 0
$ws.Range("Q74").Value = 0
$ws.Range("R74").Value = ''

# Row 75: A 19693-2025
$ws.Range("A75").Value = 'A 19693-2025'
$ws.Range("B75").Value = 45771.25997685185
$ws.Range("C75").Value = 46077
$ws.Range("D75").Value = 'SKÅNE LÄN'
$ws.Range("E75").Value = 'SIMRISHAMN'
$ws.Range("F75").Value = ''
$ws.Range("G75").Value = 1.7
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = 0
$ws.Range("N75").Value = 0
$ws.Range("O75").Value = 0
$ws.Range("P75").Value = 0
$ws.Range("Q75").Value = 0
$ws.Range("R75").Value = ''

# Row 76: A 53978-2022
$ws.Range("A76").Value = 'A 53978-2022'
$ws.Range("B76").Value = 44876
$ws.Range("C76").Value = 46077
$ws.Range("D76").Value = 'SKÅNE LÄN'
$ws.Range("E76").Value = 'SIMRISHAMN'
$ws.Range("F76").Value = ''
$ws.Range("G76").Value = 1.7
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = 0
$ws.Range("N76").Value = 0
$ws.Range("O76").Value = 0
$ws.Range("P76").Value = 0
$ws.Range("Q76").Value = 0
$ws.Range("R76").Value = ''

# Row 77: A 20096-2024
$ws.Range("A77").Value = 'A 20096-2024'
$ws.Range("B77").Value = 45434
$ws.Range("C77").Value = 46077
$ws.Range("D77").Value = 'SKÅNE LÄN'
$ws.Range("E77").Value = 'SIMRISHAMN'
$ws.Range("F77").Value = ''
$ws.Range("G77").Value = 0.8
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = 0
$ws.Range("N77").Value = 0
$ws.Range("O77").Value = 0
$ws.Range("P77").Value = 0
$ws.Range("Q77").Value = 0
$ws.Range("R77").Value = ''

# Row 78: A 3604-2024
$ws.Range("A78").Value = 'A 3604-2024'
$ws.Range("B78").Value = 45320
$ws.Range("C78").Value = 46077
$ws.Range("D78").Value = 'SKÅNE LÄN'
$ws.Range("E78").Value = 'SIMRISHAMN'
$ws.Range("F78").Value = ''
$ws.Range("G78").Value = 4.1
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = 0
$ws.Range("N78").Value = 0
$ws.Range("O78").Value = 0
$ws.Range("P78").Value = 0
$ws.Range("Q78").Value = 0
$ws.Range("R78").Value = ''

# Row 79: A 2554-2023
$ws.Range("A79").Value = 'A 2554-2023'
$ws.Range("B79").Value = 44943.67172453704
$ws.Range("C79").Value = 46077
$ws.Range("D79").Value = 'SKÅNE LÄN'
$ws.Range("E79").Value = 'SIMRISHAMN'
$ws.Range("F79").Value = ''
$ws.Range("G79").Value = 0.1
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = 0
$ws.Range("N79").Value = 0
$ws.Range("O79").Value = 0
$ws.Range("P79").Value = 0
$ws.Range("Q79").Value = 0
$ws.Range("R79").Value = ''

# Row 80: A 43851-2024
$ws.Range("A80").Value = 'A 43851-2024'
$ws.Range("B80").Value = 45572.35173611111
$ws.Range("C80").Value = 46077
$ws.Range("D80").Value = 'SKÅNE LÄN'
$ws.Range("E80").Value = 'SIMRISHAMN'
$ws.Range("F80").Value = ''
$ws.Range("G80").Value = 0.7
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = 0
$ws.Range("N80").Value = 0
$ws.Range("O80").Value = 0
$ws.Range("P80").Value = 0
$ws.Range("Q80").Value = 0
$ws.Range("R80").Value = ''

# Row 81: A 42016-2022
$ws.Range("A81").Value = 'A 42016-2022'
$ws.Range("B81").Value = 44827
$ws.Range("C81").Value = 46077
$ws.Range("D81").Value = 'SKÅNE LÄN'
$ws.Range("E81").Value = 'SIMRISHAMN'
$ws.Range("F81").Value = ''
$ws.Range("G81").Value = 0.8
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = 0
$ws.Range("N81").Value = 0
$ws.Range("O81").Value = 0
$ws.Range("P81").Value = 0
$ws.Range("Q81").Value = 0
$ws.Range("R81").Value = ''

# Row 82: A 27426-2024
$ws.Range("A82").Value = 'A 27426-2024'
$ws.Range("B82").Value = 45474.39143518519
$ws.Range("C82").Value = 46077
$ws.Range("D82").Value = 'SKÅNE LÄN'
$ws.Range("E82").Value = 'SIMRISHAMN'
$ws.Range("F82").Value = ''
$ws.Range("G82").Value = 0.6
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = 0
$ws.Range("N82").Value = 0
$ws.Range("O82").Value = 0
$ws.Range("P82").Value = 0
$ws.Range("Q82").Value = 0
$ws.Range("R82").Value = ''

# Row 83: A 45946-2022
$ws.Range("A83").Value = 'A 45946-2022'
$ws.Range("B83").Value = 44845
$ws.Range("C83").Value = 46077
$ws.Range("D83").Value = 'SKÅNE LÄN'
$ws.Range("E83").Value = 'SIMRISHAMN'
$ws.Range("F83").Value = ''
$ws.Range("G83").Value = 14
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = 0
$ws.Range("N83").Value = 0
$ws.Range("O83").Value = 0
$ws.Range("P83").Value = 0
$ws.Range("Q83").Value = 0
$ws.Range("R83").Value = ''

# Row 84: A 49303-2025
$ws.Range("A84").Value = 'A 49303-2025'
$ws.Range("B84").Value = 45938.55403935185
$ws.Range("C84").Value = 46077
$ws.Range("D84").Value = 'SKÅNE LÄN'
$ws.Range("E84").Value = 'SIMRISHAMN'
$ws.Range("F84").Value = 'Övriga Aktiebolag'
$ws.Range("G84").Value = 4.5
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = 0
$ws.Range("N84").Value = 0
$ws.Range("O84").Value = 0
$ws.Range("P84").Value = 0
$ws.Range("Q84").Value = 0
$ws.Range("R84").Value = ''

# Row 85: A 42231-2023
$ws.Range("A85").Value = 'A 42231-2023'
$ws.Range("B85").Value = 45180.31387731482
$ws.Range("C85").Value = 46077
$ws.Range("D85").Value = 'SKÅNE LÄN'
$ws.Range("E85").Value = 'SIMRISHAMN'
$ws.Range("F85").Value = ''
$ws.Range("G85").Value = 2.4
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = 0
$ws.Range("N85").Value = 0
$ws.Range("O85").Value = 0
$ws.Range("P85").Value = 0
$ws.Range("Q85").Value = 0
$ws.Range("R85").Value = ''

# Row 86: A 19137-2025
$ws.Range("A86").Value = 'A 19137-2025'
$ws.Range("B86").Value = 45769
$ws.Range("C86").Value = 46077
$ws.Range("D86").Value = 'SKÅNE LÄN'
$ws.Range("E86").Value = 'SIMRISHAMN'
$ws.Range("F86").Value = ''
$ws.Range("G86").Value = 2.2
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = 0
$ws.Range("N86").Value = 0
$ws.Range("O86").Value = 0
$ws.Range("P86").Value = 0
$ws.Range("Q86").Value = 0
$ws.Range("R86").Value = ''

# Row 87: A 19141-2025
$ws.Range("A87").Value = 'A 19141-2025'
$ws.Range("B87").Value = 45769
$ws.Range("C87").Value = 46077
$ws.Range("D87").Value = 'SKÅNE LÄN'
$ws.Range("E87").Value = 'SIMRISHAMN'
$ws.Range("F87").Value = ''
$ws.Range("G87").Value = 0.8
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = 0
$ws.Range("N87").Value = 0
$ws.Range("O87").Value = 0
$ws.Range("P87").Value = 0
$ws.Range("Q87").Value = 0
$ws.Range("R87").Value = ''

# Row 88: A 28269-2025
$ws.Range("A88").Value = 'A 28269-2025'
$ws.Range("B88").Value = 45818
$ws.Range("C88").Value = 46077
$ws.Range("D88").Value = 'SKÅNE LÄN'
$ws.Range("E88").Value = 'SIMRISHAMN'
$ws.Range("F88").Value = ''
$ws.Range("G88").Value = 6.7
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = 0
$ws.Range("N88").Value = 0
$ws.Range("O88").Value = 0
$ws.Range("P88").Value = 0
$ws.Range("Q88").Value = 0
$ws.Range("R88").Value = ''

# Row 89: A 14103-2024
$ws.Range("A89").Value = 'A 14103-2024'
$ws.Range("B89").Value = 45392
$ws.Range("C89").Value = 46077
$ws.Range("D89").Value = 'SKÅNE LÄN'
$ws.Range("E89").Value = 'SIMRISHAMN'
$ws.Range("F89").Value = ''
$ws.Range("G89").Value = 0.4
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = 0
$ws.Range("N89").Value = 0
$ws.Range("O89").Value = 0
$ws.Range("P89").Value = 0
$ws.Range("Q89").Value = 0
$ws.Range("R89").Value = ''

# Row 90: A 22194-2023
$ws.Range("A90").Value = 'A 22194-2023'
$ws.Range("B90").Value = 45069
$ws.Range("C90").Value = 46077
$ws.Range("D90").Value = 'SKÅNE LÄN'
$ws.Range("E90").Value = 'SIMRISHAMN'
$ws.Range("F90").Value = 'Övriga Aktiebolag'
$ws.Range("G90").Value = 11.4
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("M90").Value = 0
$ws.Range("N90").Value = 0
$ws.Range("O90").Value = 0
$ws.Range("P90").Value = 0
$ws.Range("Q90").Value = 0
$ws.Range("R90").Value = ''

# Row 91: A 16258-2024
$ws.Range("A91").Value = 'A 16258-2024'
$ws.Range("B91").Value = 45407
$ws.Range("C91").Value = 46077
$ws.Range("D91").Value = 'SKÅNE LÄN'
$ws.Range("E91").Value = 'SIMRISHAMN'
$ws.Range("F91").Value = ''
$ws.Range("G91").Value = 0.5
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = 0
$ws.Range("N91").Value = 0
$ws.Range("O91").Value = 0
$ws.Range("P91").Value = 0
$ws.Range("Q91").Value = 0
$ws.Range("R91").Value = ''

# Row 92: A 50138-2024
$ws.Range("A92").Value = 'A 50138-2024'
$ws.Range("B92").Value = 45600.44368055555
$ws.Range("C92").Value = 46077
$ws.Range("D92").Value = 'SKÅNE LÄN'
$ws.Range("E92").Value = 'SIMRISHAMN'
$ws.Range("F92").Value = ''
$ws.Range("G92").Value = 1.4
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 0
$ws.Range("N92").Value = 0
$ws.Range("O92").Value = 0
$ws.Range("P92").Value = 0
$ws.Range("Q92").Value = 0
$ws.Range("R92").Value = ''

# Row 93: A 55504-2025
$ws.Range("A93").Value = 'A 55504-2025'
$ws.Range("B93").Value = 45971.59123842593
$ws.Range("C93").Value = 46077
$ws.Range("D93").Value = 'SKÅNE LÄN'
$ws.Range("E93").Value = 'SIMRISHAMN'
$ws.Range("F93").Value = ''
$ws.Range("G93").Value = 0.6
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 0
$ws.Range("N93").Value = 0
$ws.Range("O93").Value = 0
$ws.Range("P93").Value = 0
$ws.Range("Q93").Value = 0
$ws.Range("R93").Value = ''

# Row 94: A 60059-2025
$ws.Range("A94").Value = 'A 60059-2025'
$ws.Range("B94").Value = 45993.65428240741
$ws.Range("C94").Value = 46077
$ws.Range("D94").Value = 'SKÅNE LÄN'
$ws.Range("E94").Value = 'SIMRISHAMN'
$ws.Range("F94").Value = 'Övriga Aktiebolag'
$ws.Range("G94").Value = 7.8
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = 0
$ws.Range("N94").Value = 0
$ws.Range("O94").Value = 0
$ws.Range("P94").Value = 0
$ws.Range("Q94").Value = 0
$ws.Range("R94").Value = ''

# Row 95: A 47874-2023
$ws.Range("A95").Value = 'A 47874-2023'
$ws.Range("B95").Value = 45204
$ws.Range("C95").Value = 46077
$ws.Range("D95").Value = 'SKÅNE LÄN'
$ws.Range("E95").Value = 'SIMRISHAMN'
$ws.Range("F95").Value = 'Övriga Aktiebolag'
$ws.Range("G95").Value = 1.4
$ws.Range("H95").Value = 0
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("M95").Value = 0
$ws.Range("N95").Value = 0
$ws.Range("O95").Value = 0
$ws.Range("P95").Value = 0
$ws.Range("Q95").Value = 0
$ws.Range("R95").Value = ''

# Row 96: A 23052-2025
$ws.Range("A96").Value = 'A 23052-2025'
$ws.Range("B96").Value = 45790.71023148148
$ws.Range("C96").Value = 46077
$ws.Range("D96").Value = 'SKÅNE LÄN'
$ws.Range("E96").Value = 'SIMRISHAMN'
$ws.Range("F96").Value = ''
$ws.Range("G96").Value = 1.2
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = 0
$ws.Range("N96").Value = 0
$ws.Range("O96").Value = 0
$ws.Range("P96").Value = 0
$ws.Range("Q96").Value = 0
$ws.Range("R96").Value = ''

# Row 97: A 12000-2025
$ws.Range("A97").Value = 'A 12000-2025'
$ws.Range("B97").Value = 45728.61288194444
$ws.Range("C97").Value = 46077
$ws.Range("D97").Value = 'SKÅNE LÄN'
$ws.Range("E97").Value = 'SIMRISHAMN'
$ws.Range("F97").Value = 'Övriga Aktiebolag'
$ws.Range("G97").Value = 2.4
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = 0
$ws.Range("N97").Value = 0
$ws.Range("O97").Value = 0
$ws.Range("P97").Value = 0
$ws.Range("Q97").Value = 0
$ws.Range("R97").Value = ''

# Row 98: A 62316-2025
$ws.Range("A98").Value = 'A 62316-2025'
$ws.Range("B98").Value = 46006.69094907407
$ws.Range("C98").Value = 46077
$ws.Range("D98").Value = 'SKÅNE LÄN'
$ws.Range("E98").Value = 'SIMRISHAMN'
$ws.Range("F98").Value = ''
$ws.Range("G98").Value = 3.5
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 0
$ws.Range("N98").Value = 0
$ws.Range("O98").Value = 0
$ws.Range("P98").Value = 0
$ws.Range("Q98").Value = 0
$ws.Range("R98").Value = ''

# Row 99: A 62357-2025
$ws.Range("A99").Value = 'A 62357-2025'
$ws.Range("B99").Value = 46007
$ws.Range("C99").Value = 46077
$ws.Range("D99").Value = 'SKÅNE LÄN'
$ws.Range("E99").Value = 'SIMRISHAMN'
$ws.Range("F99").Value = 'Övriga Aktiebolag'
$ws.Range("G99").Value = 5.5
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 0
$ws.Range("N99").Value = 0
$ws.Range("O99").Value = 0
$ws.Range("P99").Value = 0
$ws.Range("Q99").Value = 0
$ws.Range("R99").Value = ''

# Row 100: A 1910-2026
$ws.Range("A100").Value = 'A 1910-2026'
$ws.Range("B100").Value = 46035
$ws.Range("C100").Value = 46077
$ws.Range("D100").Value = 'SKÅNE LÄN'
$ws.Range("E100").Value = 'SIMRISHAMN'
$ws.Range("F100").Value = ''
$ws.Range("G100").Value = 0.7
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = 0
$ws.Range("N100").Value = 0
$ws.Range("O100").Value = 0
$ws.Range("P100").Value = 0
$ws.Range("Q100").Value = 0
$ws.Range("R100").Value = ''

# Row 101: A 30775-2025
$ws.Range("A101").Value = 'A 30775-2025'
$ws.Range("B101").Value = 45831.6150462963
$ws.Range("C101").Value = 46077
$ws.Range("D101").Value = 'SKÅNE LÄN'
$ws.Range("E101").Value = 'SIMRISHAMN'
$ws.Range("F101").Value = ''
$ws.Range("G101").Value = 0.8
$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("M101").Value = 0
$ws.Range("N101").Value = 0
$ws.Range("O101").Value = 0
$ws.Range("P101").Value = 0
$ws.Range("Q101").Value = 0
$ws.Range("R101").Value = ''

# Row 102: A 32488-2025
$ws.Range("A102").Value = 'A 32488-2025'
$ws.Range("B102").Value = 45838
$ws.Range("C102").Value = 46077
$ws.Range("D102").Value = 'SKÅNE LÄN'
$ws.Range("E102").Value = 'SIMRISHAMN'
$ws.Range("F102").Value = ''
$ws.Range("G102").Value = 4.7
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 0
$ws.Range("N102").Value = 0
$ws.Range("O102").Value = 0
$ws.Range("P102").Value = 0
$ws.Range("Q102").Value = 0
$ws.Range("R102").Value = ''

# Row 103: A 6830-2025
$ws.Range("A103").Value = 'A 6830-2025'
$ws.Range("B103").Value = 45700
$ws.Range("C103").Value = 46077
$ws.Range("D103").Value = 'SKÅNE LÄN'
$ws.Range("E103").Value = 'SIMRISHAMN'
$ws.Range("F103").Value = ''
$ws.Range("G103").Value = 1.7
$ws.Range("H103").Value = 0
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = 0
$ws.Range("N103").Value = 0
$ws.Range("O103").Value = 0
$ws.Range("P103").Value = 0
$ws.Range("Q103").Value = 0
$ws.Range("R103").Value = ''

# Row 104: A 37189-2025
$ws.Range("A104").Value = 'A 37189-2025'
$ws.Range("B104").Value = 45875
$ws.Range("C104").Value = 46077
$ws.Range("D104").Value = 'SKÅNE LÄN'
$ws.Range("E104").Value = 'SIMRISHAMN'
$ws.Range("F104").Value = ''
$ws.Range("G104").Value = 1.2
$ws.Range("H104").Value = 0
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("M104").Value = 0
$ws.Range("N104").Value = 0
$ws.Range("O104").Value = 0
$ws.Range("P104").Value = 0
$ws.Range("Q104").Value = 0
$ws.Range("R104").Value = ''

# Row 105: A 37052-2025
$ws.Range("A105").Value = 'A 37052-2025'
$ws.Range("B105").Value = 45875
$ws.Range("C105").Value = 46077
$ws.Range("D105").Value = 'SKÅNE LÄN'
$ws.Range("E105").Value = 'SIMRISHAMN'
$ws.Range("F105").Value = ''
$ws.Range("G105").Value = 0.9
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 0
$ws.Range("N105").Value = 0
$ws.Range("O105").Value = 0
$ws.Range("P105").Value = 0
$ws.Range("Q105").Value = 0
$ws.Range("R105").Value = ''

# Fix row heights (avoid autofit side effects from multi-line text)
for ($r = 2; $r -le 105; $r++) {
    $ws.Rows.Item($r).RowHeight = 15
}